$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.446407857441102
$ws.Range("C2").Value = 0.04875525018707094
$ws.Range("D2").Value = 0.1248802924412242
$ws.Range("F2").Value = 2.319677052181873
$ws.Range("G2").Value = 0.002555383951192388
$ws.Range("K2").Value = 0.9248430034724322
$ws.Range("L2").Value = 0.3072403409413056
$ws.Range("N2").Value = 3.149778894070607

$ws.Range("B3").Value = 1.393681655594889
$ws.Range("C3").Value = 0.04570271063329301
$ws.Range("D3").Value = 0.1252569564732422
$ws.Range("F3").Value = 2.292421660901397
$ws.Range("G3").Value = 0.002559777865054426
$ws.Range("K3").Value = 0.8758793560311631
$ws.Range("L3").Value = 0.2967189422946461
$ws.Range("N3").Value = 3.150523909750476

$ws.Range("B4").Value = 1.3621963007528
$ws.Range("C4").Value = 0.04380598120572898
$ws.Range("D4").Value = 0.1254937250592345
$ws.Range("F4").Value = 2.276716048482683
$ws.Range("G4").Value = 0.002562618322567956
$ws.Range("K4").Value = 0.8463997843549862
$ws.Range("L4").Value = 0.2904453072818853
$ws.Range("N4").Value = 3.15169223227474

$ws.Range("B5").Value = 1.349588914503272
$ws.Range("C5").Value = 0.04302734922904961
$ws.Range("D5").Value = 0.1255915916881003
$ws.Range("F5").Value = 2.27057424750241
$ws.Range("G5").Value = 0.002563811801801901
$ws.Range("K5").Value = 0.8345331405825505
$ws.Range("L5").Value = 0.2879355311150675
$ws.Range("N5").Value = 3.152346779376458

$ws.Range("B6").Value = 1.347508936042573
$ws.Range("C6").Value = 0.04289771223904637
$ws.Range("D6").Value = 0.1256079257999758
$ws.Range("F6").Value = 2.269569997272086
$ws.Range("G6").Value = 0.002564012154224744
$ws.Range("K6").Value = 0.8325715312619195
$ws.Range("L6").Value = 0.2875216064793875
$ws.Range("N6").Value = 3.152466235131286

$ws.Range("B7").Value = 1.362025370034075
$ws.Range("C7").Value = 0.04379550344730632
$ws.Range("D7").Value = 0.1254950393298255
$ws.Range("F7").Value = 2.276632172527343
$ws.Range("G7").Value = 0.002562634272524848
$ws.Range("K7").Value = 0.8462391539221983
$ws.Range("L7").Value = 0.2904112703134132
$ws.Range("N7").Value = 3.151700337636854

$ws.Range("B8").Value = 1.428043181316554
$ws.Range("C8").Value = 0.04770736349097859
$ws.Range("D8").Value = 0.1250090254620337
$ws.Range("F8").Value = 2.310065465743264
$ws.Range("G8").Value = 0.002556869451732701
$ws.Range("K8").Value = 0.9078388002073439
$ws.Range("L8").Value = 0.3035737637162867
$ws.Range("N8").Value = 3.149888041283745

$ws.Range("B9").Value = 1.564579758614514
$ws.Range("C9").Value = 0.05520296890088616
$ws.Range("D9").Value = 0.1240995486631338
$ws.Range("F9").Value = 2.383822121437589
$ws.Range("G9").Value = 0.002546690563854503
$ws.Range("K9").Value = 1.033297824626231
$ws.Range("L9").Value = 0.3308730733397738
$ws.Range("N9").Value = 3.151990310222658

$ws.Range("B10").Value = 1.66925326720542
$ws.Range("C10").Value = 0.06060708034839024
$ws.Range("D10").Value = 0.1234578875293941
$ws.Range("F10").Value = 2.443052514138344
$ws.Range("G10").Value = 0.002539890986359187
$ws.Range("K10").Value = 1.128362118645725
$ws.Range("L10").Value = 0.3518498333712472
$ws.Range("N10").Value = 3.157007100923181

$ws.Range("B11").Value = 1.717829692625628
$ws.Range("C11").Value = 0.06304411620199346
$ws.Range("D11").Value = 0.1231717322876413
$ws.Range("F11").Value = 2.471103236394839
$ws.Range("G11").Value = 0.002536943482880707
$ws.Range("K11").Value = 1.172247760886194
$ws.Range("L11").Value = 0.3615954565542978
$ws.Range("N11").Value = 3.160048765877434

$ws.Range("B12").Value = 1.736362993197474
$ws.Range("C12").Value = 0.06396395630953577
$ws.Range("D12").Value = 0.1230641989681853
$ws.Range("F12").Value = 2.481885131951628
$ws.Range("G12").Value = 0.002535848163094223
$ws.Range("K12").Value = 1.188958874148682
$ws.Range("L12").Value = 0.3653152811983631
$ws.Range("N12").Value = 3.161310191307933

$ws.Range("B13").Value = 1.732365350225109
$ws.Range("C13").Value = 0.06376598540882128
$ws.Range("D13").Value = 0.1230873213557135
$ws.Range("F13").Value = 2.479555945977154
$ws.Range("G13").Value = 0.002536083134889926
$ws.Range("K13").Value = 1.185355711919414
$ws.Range("L13").Value = 0.3645128418089172
$ws.Range("N13").Value = 3.161033639312137

$ws.Range("B14").Value = 1.719351660932602
$ws.Range("C14").Value = 0.06311985206622239
$ws.Range("D14").Value = 0.1231628688606907
$ws.Range("F14").Value = 2.471987065535018
$ws.Range("G14").Value = 0.002536852953155305
$ws.Range("K14").Value = 1.173620736135291
$ws.Range("L14").Value = 0.3619008993056951
$ws.Range("N14").Value = 3.160150345005405

$ws.Range("B15").Value = 1.711398445270106
$ws.Range("C15").Value = 0.06272368601990763
$ws.Range("D15").Value = 0.1232092517148669
$ws.Range("F15").Value = 2.467371722195338
$ws.Range("G15").Value = 0.002537327199389964
$ws.Range("K15").Value = 1.166444794834121
$ws.Range("L15").Value = 0.3603048382693004
$ws.Range("N15").Value = 3.159623588464612

$ws.Range("B16").Value = 1.666098035851121
$ws.Range("C16").Value = 0.06044738984809328
$ws.Range("D16").Value = 0.1234767041189571
$ws.Range("F16").Value = 2.441241648799362
$ws.Range("G16").Value = 0.002540086534323181
$ws.Range("K16").Value = 1.125507016504088
$ws.Range("L16").Value = 0.3512170354602802
$ws.Range("N16").Value = 3.156823637933854

$ws.Range("B17").Value = 1.63855386072197
$ws.Range("C17").Value = 0.05904553288958425
$ws.Range("D17").Value = 0.1236422489952016
$ws.Range("F17").Value = 2.425495510664049
$ws.Range("G17").Value = 0.002541816525370155
$ws.Range("K17").Value = 1.100557381378877
$ws.Range("L17").Value = 0.3456941106902036
$ws.Range("N17").Value = 3.155300766647116

$ws.Range("B18").Value = 1.622801505766915
$ws.Range("C18").Value = 0.05823721173537422
$ws.Range("D18").Value = 0.1237380061312496
$ws.Range("F18").Value = 2.416542864107754
$ws.Range("G18").Value = 0.002542825287083378
$ws.Range("K18").Value = 1.086267260912109
$ws.Range("L18").Value = 0.3425365865015522
$ws.Range("N18").Value = 3.154496315176218

$ws.Range("B19").Value = 1.61748352981715
$ws.Range("C19").Value = 0.05796318121892341
$ws.Range("D19").Value = 0.1237705205536717
$ws.Range("F19").Value = 2.413529516255011
$ws.Range("G19").Value = 0.002543169195822204
$ws.Range("K19").Value = 1.081439205103692
$ws.Range("L19").Value = 0.3414707817300808
$ws.Range("N19").Value = 3.154236203672653

$ws.Range("B20").Value = 1.641476633738193
$ws.Range("C20").Value = 0.05919497060627066
$ws.Range("D20").Value = 0.1236245705487509
$ws.Range("F20").Value = 2.427160934242721
$ws.Range("G20").Value = 0.002541630946211064
$ws.Range("K20").Value = 1.103207073554785
$ws.Range("L20").Value = 0.346280055937342
$ws.Range("N20").Value = 3.155455479845102

$ws.Range("B21").Value = 1.723170335636496
$ws.Range("C21").Value = 0.06330971842109534
$ws.Range("D21").Value = 0.1231406562521222
$ws.Range("F21").Value = 2.474205890425509
$ws.Range("G21").Value = 0.002536626274308514
$ws.Range("K21").Value = 1.17706506668921
$ws.Range("L21").Value = 0.3626672918673393
$ws.Range("N21").Value = 3.160406811894603

$ws.Range("B22").Value = 1.77736916496076
$ws.Range("C22").Value = 0.06598141006804781
$ws.Range("D22").Value = 0.122829214489391
$ws.Range("F22").Value = 2.505883594346244
$ws.Range("G22").Value = 0.002533476829619233
$ws.Range("K22").Value = 1.225875351910389
$ws.Range("L22").Value = 0.3735485523019122
$ws.Range("N22").Value = 3.164281890062853

$ws.Range("B23").Value = 1.748368219892711
$ws.Range("C23").Value = 0.06455706452130983
$ws.Range("D23").Value = 0.1229949946533324
$ws.Range("F23").Value = 2.488891217314716
$ws.Range("G23").Value = 0.002535146675549111
$ws.Range("K23").Value = 1.199774828943276
$ws.Range("L23").Value = 0.3677252987585433
$ws.Range("N23").Value = 3.162155077804641

$ws.Range("B24").Value = 1.640154989065309
$ws.Range("C24").Value = 0.05912741721638781
$ws.Range("D24").Value = 0.1236325611574909
$ws.Range("F24").Value = 2.426407684735722
$ws.Range("G24").Value = 0.002541714802235974
$ws.Range("K24").Value = 1.102008980315702
$ws.Range("L24").Value = 0.346015095038041
$ws.Range("N24").Value = 3.155385312655383

$ws.Range("B25").Value = 1.52688072056867
$ws.Range("C25").Value = 0.0531935788532536
$ws.Range("D25").Value = 0.1243409245502765
$ws.Range("F25").Value = 2.362987034436017
$ws.Range("G25").Value = 0.002549324467760788
$ws.Range("K25").Value = 0.9988538513990477
$ws.Range("L25").Value = 0.3233272813683072
$ws.Range("N25").Value = 3.150813587003128
